$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2096774193548387
$ws.Range("C2").Value = 0.5204301075268817
$ws.Range("J2").Value = 0.01720430107526882
$ws.Range("P2").Value = 0.1451612903225807
$ws.Range("S2").Value = 0.1075268817204301
$ws.Range("B3").Value = 0.008048289738430584
$ws.Range("C3").Value = 0.03420523138832998
$ws.Range("J3").Value = 0.02414486921529175
$ws.Range("P3").Value = 0.6921529175050302
$ws.Range("S3").Value = 0.2414486921529175
$ws.Range("J4").Value = 0.09734513274336283
$ws.Range("P4").Value = 0.6814159292035398
$ws.Range("S4").Value = 0.2212389380530974
$ws.Range("O5").Value = 0.1111111111111111
$ws.Range("P5").Value = 0.5555555555555556
$ws.Range("S5").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.07258064516129033
$ws.Range("D6").Value = 0.01451612903225807
$ws.Range("E6").Value = 0.001612903225806452
$ws.Range("F6").Value = 0.06290322580645161
$ws.Range("J6").Value = 0.232258064516129
$ws.Range("O6").Value = 0.02903225806451613
$ws.Range("Q6").Value = 0.1290322580645161
$ws.Range("R6").Value = 0.07258064516129033
$ws.Range("S6").Value = 0.3854838709677419
$ws.Range("B7").Value = 0.08904109589041095
$ws.Range("D7").Value = 0.0136986301369863
$ws.Range("F7").Value = 0.04794520547945205
$ws.Range("J7").Value = 0.1404109589041096
$ws.Range("O7").Value = 0.02054794520547945
$ws.Range("Q7").Value = 0.160958904109589
$ws.Range("R7").Value = 0.0684931506849315
$ws.Range("S7").Value = 0.4589041095890411
$ws.Range("B8").Value = 0.1029306647605432
$ws.Range("D8").Value = 0.01286633309506791
$ws.Range("E8").Value = 0.0007147962830593281
$ws.Range("F8").Value = 0.05289492494639028
$ws.Range("J8").Value = 0.1100786275911365
$ws.Range("O8").Value = 0.02144388849177984
$ws.Range("Q8").Value = 0.1672623302358828
$ws.Range("R8").Value = 0.1007862759113653
$ws.Range("S8").Value = 0.4310221586847748
$ws.Range("B9").Value = 0.09210526315789473
$ws.Range("D9").Value = 0.01503759398496241
$ws.Range("E9").Value = 0.003759398496240601
$ws.Range("F9").Value = 0.06390977443609022
$ws.Range("J9").Value = 0.112781954887218
$ws.Range("O9").Value = 0.01879699248120301
$ws.Range("Q9").Value = 0.1578947368421053
$ws.Range("R9").Value = 0.07142857142857142
$ws.Range("S9").Value = 0.4642857142857143
$ws.Range("B10").Value = 0.1124675324675325
$ws.Range("D10").Value = 0.02051948051948052
$ws.Range("E10").Value = 0.001818181818181818
$ws.Range("F10").Value = 0.06805194805194806
$ws.Range("J10").Value = 0.1207792207792208
$ws.Range("O10").Value = 0.01428571428571429
$ws.Range("Q10").Value = 0.2187012987012987
$ws.Range("R10").Value = 0.07350649350649351
$ws.Range("S10").Value = 0.3698701298701298
$ws.Range("F11").Value = 0.001020408163265306
$ws.Range("G11").Value = 0.1571428571428571
$ws.Range("J11").Value = 0.09285714285714286
$ws.Range("K11").Value = 0.2163265306122449
$ws.Range("L11").Value = 0.5163265306122449
$ws.Range("S11").Value = 0.0163265306122449
$ws.Range("G12").Value = 0.7126654064272212
$ws.Range("J12").Value = 0.2268431001890359
$ws.Range("K12").Value = 0.007561436672967864
$ws.Range("L12").Value = 0.02079395085066163
$ws.Range("S12").Value = 0.03213610586011342
$ws.Range("F13").Value = 0.007518796992481203
$ws.Range("G13").Value = 0.5864661654135338
$ws.Range("J13").Value = 0.3383458646616541
$ws.Range("S13").Value = 0.06766917293233082
$ws.Range("F15").Value = 0.01567398119122257
$ws.Range("H15").Value = 0.1457680250783699
$ws.Range("I15").Value = 0.04858934169278997
$ws.Range("J15").Value = 0.3699059561128527
$ws.Range("K15").Value = 0.06583072100313479
$ws.Range("M15").Value = 0.0109717868338558
$ws.Range("O15").Value = 0.054858934169279
$ws.Range("S15").Value = 0.2884012539184953
$ws.Range("F16").Value = 0.01470588235294118
$ws.Range("H16").Value = 0.1819852941176471
$ws.Range("I16").Value = 0.08272058823529412
$ws.Range("J16").Value = 0.3841911764705883
$ws.Range("K16").Value = 0.1066176470588235
$ws.Range("M16").Value = 0.02573529411764706
$ws.Range("O16").Value = 0.04227941176470588
$ws.Range("S16").Value = 0.1617647058823529
$ws.Range("F17").Value = 0.01355421686746988
$ws.Range("H17").Value = 0.1626506024096386
$ws.Range("I17").Value = 0.08207831325301204
$ws.Range("J17").Value = 0.427710843373494
$ws.Range("K17").Value = 0.1031626506024096
$ws.Range("M17").Value = 0.01506024096385542
$ws.Range("O17").Value = 0.0572289156626506
$ws.Range("S17").Value = 0.1385542168674699
$ws.Range("F18").Value = 0.02205882352941177
$ws.Range("H18").Value = 0.1452205882352941
$ws.Range("I18").Value = 0.08455882352941177
$ws.Range("J18").Value = 0.4172794117647059
$ws.Range("K18").Value = 0.1194852941176471
$ws.Range("M18").Value = 0.01102941176470588
$ws.Range("O18").Value = 0.05330882352941176
$ws.Range("S18").Value = 0.1470588235294118
$ws.Range("F19").Value = 0.01555285540704739
$ws.Range("H19").Value = 0.2240583232077764
$ws.Range("I19").Value = 0.07363304981773998
$ws.Range("J19").Value = 0.3538274605103281
$ws.Range("K19").Value = 0.1105710814094775
$ws.Range("M19").Value = 0.02162818955042527
$ws.Range("N19").Value = 0.0009720534629404617
$ws.Range("O19").Value = 0.06366950182260024
$ws.Range("S19").Value = 0.1360874848116646
